$wb = $excel.ActiveWorkbook
$regWs = $wb.Worksheets.Item("registration")
$loginWs = $wb.Worksheets.Item("login")

# Copy formatting (styles/borders) from row 19 down into rows 20-28
$regWs.Range("A19:P19").Copy()
$regWs.Range("A20:P28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column C
$regWs.Range("C20").Value = "Khanb"
$regWs.Range("C21").Value = "Khanc"
$regWs.Range("C22").Value = "Khand"
$regWs.Range("C23").Value = "Khane"
$regWs.Range("C24").Value = "Khanf"
$regWs.Range("C25").Value = "Khang"
$regWs.Range("C26").Value = "Khanh"
$regWs.Range("C27").Value = "Khani"
$regWs.Range("C28").Value = "Khanj"

# Column D
$regWs.Range("D20").Value = "fname19@gmail.com"
$regWs.Range("D21").Value = "fname20@gmail.com"
$regWs.Range("D22").Value = "fname21@gmail.com"
$regWs.Range("D23").Value = "fname22@gmail.com"
$regWs.Range("D24").Value = "fname23@gmail.com"
$regWs.Range("D25").Value = "fname24@gmail.com"
$regWs.Range("D26").Value = "fname25@gmail.com"
$regWs.Range("D27").Value = "fname26@gmail.com"
$regWs.Range("D28").Value = "fname27@gmail.com"

# Column E
$regWs.Range("E20").Value = "23/5/1995"
$regWs.Range("E21").Value = "23/5/1996"
$regWs.Range("E22").Value = "23/5/1997"
$regWs.Range("E23").Value = "23/5/1998"
$regWs.Range("E24").Value = "23/5/1999"
$regWs.Range("E25").Value = "23/5/2000"
$regWs.Range("E26").Value = "23/5/2001"
$regWs.Range("E27").Value = "23/5/2002"
$regWs.Range("E28").Value = "23/5/2003"

# Column F
$regWs.Range("F20").Value = "test@1234585"
$regWs.Range("F21").Value = "test@1234586"
$regWs.Range("F22").Value = "test@1234587"
$regWs.Range("F23").Value = "test@1234588"
$regWs.Range("F24").Value = "test@1234589"
$regWs.Range("F25").Value = "test@1234590"
$regWs.Range("F26").Value = "test@1234591"
$regWs.Range("F27").Value = "test@1234592"
$regWs.Range("F28").Value = "test@1234593"

# Column G
$regWs.Range("G20").Value = "address1,po1,lstreet19"
$regWs.Range("G21").Value = "address1,po1,lstreet20"
$regWs.Range("G22").Value = "address1,po1,lstreet21"
$regWs.Range("G23").Value = "address1,po1,lstreet22"
$regWs.Range("G24").Value = "address1,po1,lstreet23"
$regWs.Range("G25").Value = "address1,po1,lstreet24"
$regWs.Range("G26").Value = "address1,po1,lstreet25"
$regWs.Range("G27").Value = "address1,po1,lstreet26"
$regWs.Range("G28").Value = "address1,po1,lstreet27"

# Column H
$regWs.Range("H20").Value = "Address30"
$regWs.Range("H21").Value = "Address31"
$regWs.Range("H22").Value = "Address32"
$regWs.Range("H23").Value = "Address33"
$regWs.Range("H24").Value = "Address34"
$regWs.Range("H25").Value = "Address35"
$regWs.Range("H26").Value = "Address36"
$regWs.Range("H27").Value = "Address37"
$regWs.Range("H28").Value = "Address38"

# Column M
$regWs.Range("M20").Value = "add info19"
$regWs.Range("M21").Value = "add info20"
$regWs.Range("M22").Value = "add info21"
$regWs.Range("M23").Value = "add info22"
$regWs.Range("M24").Value = "add info23"
$regWs.Range("M25").Value = "add info24"
$regWs.Range("M26").Value = "add info25"
$regWs.Range("M27").Value = "add info26"
$regWs.Range("M28").Value = "add info27"

# Column O
$regWs.Range("O20").Value = "addressalias19"
$regWs.Range("O21").Value = "addressalias20"
$regWs.Range("O22").Value = "addressalias21"
$regWs.Range("O23").Value = "addressalias22"
$regWs.Range("O24").Value = "addressalias23"
$regWs.Range("O25").Value = "addressalias24"
$regWs.Range("O26").Value = "addressalias25"
$regWs.Range("O27").Value = "addressalias26"
$regWs.Range("O28").Value = "addressalias27"

# Column B
$regWs.Range("B20").Value = "Debb"
$regWs.Range("B21").Value = "Debc"
$regWs.Range("B22").Value = "Debd"
$regWs.Range("B23").Value = "Debe"
$regWs.Range("B24").Value = "Debf"
$regWs.Range("B25").Value = "Debg"
$regWs.Range("B26").Value = "Debh"
$regWs.Range("B27").Value = "Debi"
$regWs.Range("B28").Value = "Debj"

# Column A
$regWs.Range("A20").Value = "Mr."
$regWs.Range("A21").Value = "Mrs."
$regWs.Range("A22").Value = "Mr."
$regWs.Range("A23").Value = "Mr."
$regWs.Range("A24").Value = "Mr."
$regWs.Range("A25").Value = "Mr."
$regWs.Range("A26").Value = "Mr."
$regWs.Range("A27").Value = "Mr."
$regWs.Range("A28").Value = "Mr."

# Column I
$regWs.Range("I20").Value = "Bengaluru"
$regWs.Range("I21").Value = "Bengaluru"
$regWs.Range("I22").Value = "Bengaluru"
$regWs.Range("I23").Value = "Bengaluru"
$regWs.Range("I24").Value = "Bengaluru"
$regWs.Range("I25").Value = "Bengaluru"
$regWs.Range("I26").Value = "Bengaluru"
$regWs.Range("I27").Value = "Bengaluru"
$regWs.Range("I28").Value = "Bengaluru"

# Column J
$regWs.Range("J20").Value = "New York"
$regWs.Range("J21").Value = "New York"
$regWs.Range("J22").Value = "New York"
$regWs.Range("J23").Value = "New York"
$regWs.Range("J24").Value = "New York"
$regWs.Range("J25").Value = "New York"
$regWs.Range("J26").Value = "New York"
$regWs.Range("J27").Value = "New York"
$regWs.Range("J28").Value = "New York"

# Column K
$regWs.Range("K20").Value = 12363
$regWs.Range("K21").Value = 12364
$regWs.Range("K22").Value = 12365
$regWs.Range("K23").Value = 12366
$regWs.Range("K24").Value = 12367
$regWs.Range("K25").Value = 12368
$regWs.Range("K26").Value = 12369
$regWs.Range("K27").Value = 12370
$regWs.Range("K28").Value = 12371

# Column L
$regWs.Range("L20").Value = "United States"
$regWs.Range("L21").Value = "United States"
$regWs.Range("L22").Value = "United States"
$regWs.Range("L23").Value = "United States"
$regWs.Range("L24").Value = "United States"
$regWs.Range("L25").Value = "United States"
$regWs.Range("L26").Value = "United States"
$regWs.Range("L27").Value = "United States"
$regWs.Range("L28").Value = "United States"

# Column N
$regWs.Range("N20").Value = 123456807
$regWs.Range("N21").Value = 123456808
$regWs.Range("N22").Value = 123456809
$regWs.Range("N23").Value = 123456810
$regWs.Range("N24").Value = 123456811
$regWs.Range("N25").Value = 123456812
$regWs.Range("N26").Value = 123456813
$regWs.Range("N27").Value = 123456814
$regWs.Range("N28").Value = 123456815

# Column P
$regWs.Range("P20").Value = "Used"
$regWs.Range("P21").Value = "Used"
$regWs.Range("P22").Value = "Used"

# Hyperlinks for column D (email -> mailto)
$regWs.Hyperlinks.Add($regWs.Range("D20"), "mailto:fname19@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D21"), "mailto:fname20@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D22"), "mailto:fname21@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D23"), "mailto:fname22@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D24"), "mailto:fname23@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D25"), "mailto:fname24@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D26"), "mailto:fname25@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D27"), "mailto:fname26@gmail.com")
$regWs.Hyperlinks.Add($regWs.Range("D28"), "mailto:fname27@gmail.com")

# Hyperlinks for column F (password -> mailto)
$regWs.Hyperlinks.Add($regWs.Range("F20"), "mailto:test@1234585")
$regWs.Hyperlinks.Add($regWs.Range("F21:F28"), "mailto:test@1234586")

# Fill previously blank status cells
$regWs.Range("P18").Value = "Used"
$regWs.Range("P19").Value = "Used"

# login sheet: add row 5 referencing the newly-registered account
$loginWs.Range("A5").Value = "fname21@gmail.com"
$loginWs.Range("B5").Value = "test@1234587"

# Sheet view / selection updates
$regWs.Activate()
$appWindow = $excel.ActiveWindow
$appWindow.ScrollRow = 15
$appWindow.ScrollColumn = 1
$regWs.Range("B28").Select()
